$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '30.197.78'
$ws.Range('E2').Value = '  -0.28%  '

$ws.Range('D3').Value = '1.862.51'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.67'
$ws.Range('E5').Value = '  +3.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4713'
$ws.Range('E7').Value = '  +0.52%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.80'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2855'
$ws.Range('E9').Value = '  -0.25%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06472'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.73'
$ws.Range('E11').Value = '  -4.53%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07693'
$ws.Range('E12').Value = '  -3.33%  '

$ws.Range('D13').Value = '1.867.91'
$ws.Range('E13').Value = '  -0.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.27'
$ws.Range('E14').Value = '  -2.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6820'
$ws.Range('E15').Value = '  -1.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.067'
$ws.Range('E16').Value = '  -0.69%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '268.75'
$ws.Range('E17').Value = '  +0.24%  '

$ws.Range('D18').Value = '30.189.99'
$ws.Range('E18').Value = '  -0.48%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.33'
$ws.Range('E19').Value = '  -5.68%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007540'
$ws.Range('E20').Value = '  -2.86%  '

$ws.Range('E21').Value = '  -0.10%  '

$ws.Range('D22').Value = '2.109.51'
$ws.Range('E22').Value = '  -0.41%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.186'
$ws.Range('E24').Value = '  -1.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.104'
$ws.Range('E25').Value = '  -1.64%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.320'
$ws.Range('E26').Value = '  -0.55%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.50'
$ws.Range('E27').Value = '  -1.24%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.75'
$ws.Range('E28').Value = '  -0.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.373'
$ws.Range('E30').Value = '  +0.69%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09816'
$ws.Range('E31').Value = '  -0.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.516'
$ws.Range('E32').Value = '  +4.06%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.236'
$ws.Range('E33').Value = '  -2.22%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.973'
$ws.Range('E34').Value = '  -1.74%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04700'
$ws.Range('E35').Value = '  -0.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.109'
$ws.Range('E36').Value = '  -2.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6858'
$ws.Range('E37').Value = '  -2.20%  '

$ws.Range('E38').Value = '  -0.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01849'
$ws.Range('E39').Value = '  -1.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.727'
$ws.Range('E40').Value = '  -3.09%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.408'
$ws.Range('E41').Value = '  +2.65%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.30'
$ws.Range('E42').Value = '  -2.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8356'
$ws.Range('E44').Value = '  -0.67%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.883'
$ws.Range('E45').Value = '  -3.64%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.94'
$ws.Range('E46').Value = '  -0.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4062'
$ws.Range('E47').Value = '  -2.52%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.168'
$ws.Range('E48').Value = '  +0.31%  '

$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.929'
$ws.Range('E49').Value = '  -2.07%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '920.78'
$ws.Range('E50').Value = '  +0.33%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.41'
